# Insert a new weekly price record at row 56 of the "Chirimoya" sheet.
# This pushes the existing rows 56..120 down to 57..121 (dimension grows
# from A1:T120 to A1:T121), and fills the freshly inserted row 56 with the
# new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 56 (shifts 56:120 -> 57:121,
# copying the formatting -- including the date number-format on column D --
# from the row above, same as Excel's native "Insert Sheet Rows").
$ws.Rows(56).Insert()

$ws.Range("A56").Value = 10
$ws.Range("B56").Value = 'Vega Modelo de Temuco'
$ws.Range("C56").Value = 'La Araucanía'
$ws.Range("D56").Value = 44763
$ws.Range("E56").Value = 9
$ws.Range("F56").Value = 'Fruta'
$ws.Range("G56").Value = 100107
$ws.Range("H56").Value = 'Otros'
$ws.Range("I56").Value = 100107002
$ws.Range("J56").Value = 'Chirimoya'
$ws.Range("K56").Value = 'Cultivar IV Región'
$ws.Range("L56").Value = 'Primera'
$ws.Range("M56").Value = 100
$ws.Range("N56").Value = 3500
$ws.Range("O56").Value = 3500
$ws.Range("P56").Value = 3500
$ws.Range("Q56").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R56").Value = 'Provincia del Elquí'
$ws.Range("S56").Value = 3500
$ws.Range("T56").Value = 1
